$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-09 18:48:32"
$ws.Range("E3").Value = "2026-02-09 18:48:34"
$ws.Range("O3").Value = "-3.5 °C"
$ws.Range("E4").Value = "2026-02-09 18:48:36"
$ws.Range("E5").Value = "2026-02-09 18:48:39"
$ws.Range("E6").Value = "2026-02-09 18:48:41"
$ws.Range("E7").Value = "2026-02-09 18:48:43"
$ws.Range("O7").Value = "12.7 °C"
$ws.Range("E8").Value = "2026-02-09 18:48:46"
$ws.Range("E9").Value = "2026-02-09 18:48:48"
$ws.Range("E10").Value = "2026-02-09 18:48:50"
$ws.Range("O10").Value = "8.3 °C"
$ws.Range("E11").Value = "2026-02-09 18:48:53"
$ws.Range("E12").Value = "2026-02-09 18:48:55"
$ws.Range("E13").Value = "2026-02-09 18:48:57"
$ws.Range("E14").Value = "2026-02-09 18:49:00"
$ws.Range("E15").Value = "2026-02-09 18:49:02"
$ws.Range("E16").Value = "2026-02-09 18:49:04"
$ws.Range("E17").Value = "2026-02-09 18:49:07"
$ws.Range("H17").Value = "'83%"
$ws.Range("E18").Value = "2026-02-09 18:49:09"
$ws.Range("E19").Value = "2026-02-09 18:49:12"
$ws.Range("E20").Value = "2026-02-09 18:49:14"
$ws.Range("E21").Value = "2026-02-09 18:49:16"
$ws.Range("O21").Value = "4.5 °C"
$ws.Range("E22").Value = "2026-02-09 18:49:18"
$ws.Range("E23").Value = "2026-02-09 18:49:21"
$ws.Range("H23").Value = "'87%"
$ws.Range("E24").Value = "2026-02-09 18:49:23"
$ws.Range("H24").Value = "'81%"
$ws.Range("I24").Value = "0.4 mm"
$ws.Range("O24").Value = "7.9 °C"
$ws.Range("E25").Value = "2026-02-09 18:49:26"
$ws.Range("E26").Value = "2026-02-09 18:49:29"
$ws.Range("O26").Value = "2.4 °C"
$ws.Range("E27").Value = "2026-02-09 18:49:31"
$ws.Range("O27").Value = "-2.5 °C"
$ws.Range("E28").Value = "2026-02-09 18:49:33"
$ws.Range("E29").Value = "2026-02-09 18:49:36"
$ws.Range("H29").Value = "'83%"
$ws.Range("O29").Value = "8.7 °C"
$ws.Range("E30").Value = "2026-02-09 18:49:38"
$ws.Range("E31").Value = "2026-02-09 18:49:40"
$ws.Range("E32").Value = "2026-02-09 18:49:43"
$ws.Range("E33").Value = "2026-02-09 18:49:45"
$ws.Range("J33").Value = "1007.4 hPa"
$ws.Range("O33").Value = "2.9 °C"
$ws.Range("E34").Value = "2026-02-09 18:49:48"
$ws.Range("H34").Value = "'74%"
$ws.Range("E35").Value = "2026-02-09 18:49:51"
$ws.Range("I35").Value = "0.4 mm"
$ws.Range("E36").Value = "2026-02-09 18:49:53"
$ws.Range("J36").Value = "1007.2 hPa"
$ws.Range("E37").Value = "2026-02-09 18:49:56"
$ws.Range("O37").Value = "5.6 °C"
$ws.Range("E38").Value = "2026-02-09 18:49:58"
$ws.Range("E39").Value = "2026-02-09 18:50:01"
$ws.Range("E40").Value = "2026-02-09 18:50:04"
$ws.Range("O40").Value = "4.6 °C"
$ws.Range("E41").Value = "2026-02-09 18:50:06"
$ws.Range("H41").Value = "'55%"
$ws.Range("E42").Value = "2026-02-09 18:50:09"
$ws.Range("H42").Value = "'84%"
$ws.Range("E43").Value = "2026-02-09 18:50:11"
$ws.Range("E44").Value = "2026-02-09 18:50:14"
$ws.Range("O44").Value = "-4.0 °C"
$ws.Range("E45").Value = "2026-02-09 18:50:17"
$ws.Range("J45").Value = "1007.3 hPa"
$ws.Range("O45").Value = "3.8 °C"
$ws.Range("E46").Value = "2026-02-09 18:50:19"
$ws.Range("I46").Value = "0.4 mm"
